# Update cryptos Price (D) / Volume(1h) (E) columns per latest scrape.
# A leading apostrophe forces text interpretation so numeric-looking
# strings (e.g. "315.72") are stored verbatim instead of being coerced
# into real numbers by Excel's type inference (which would also drop
# meaningful trailing zeros, e.g. "0.3710" -> 0.371). The Style reset
# afterwards clears the quote-prefix formatting flag so the cell's
# style stays identical to the original (unstyled) cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.440.32"
$ws.Range("E2").Value = "'  +1.96%  "
$ws.Range("D2:E2").Style = "Normal"
$ws.Range("D3").Value = "'1.843.00"
$ws.Range("E3").Value = "'  +1.48%  "
$ws.Range("D3:E3").Style = "Normal"
$ws.Range("E4").Value = "'  +1.39%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'315.72"
$ws.Range("E5").Value = "'  +2.13%  "
$ws.Range("D5:E5").Style = "Normal"
$ws.Range("D6").Value = "'1.014"
$ws.Range("E6").Value = "'  +1.25%  "
$ws.Range("D6:E6").Style = "Normal"
$ws.Range("D7").Value = "'0.4773"
$ws.Range("E7").Value = "'  +1.80%  "
$ws.Range("D7:E7").Style = "Normal"
$ws.Range("D8").Value = "'0.3710"
$ws.Range("E8").Value = "'  +0.43%  "
$ws.Range("D8:E8").Style = "Normal"
$ws.Range("D9").Value = "'0.07475"
$ws.Range("E9").Value = "'  +1.29%  "
$ws.Range("D9:E9").Style = "Normal"
$ws.Range("D10").Value = "'0.8879"
$ws.Range("E10").Value = "'  +1.91%  "
$ws.Range("D10:E10").Style = "Normal"
$ws.Range("D11").Value = "'20.50"
$ws.Range("E11").Value = "'  +0.55%  "
$ws.Range("D11:E11").Style = "Normal"
$ws.Range("D12").Value = "'1.852.35"
$ws.Range("E12").Value = "'  +1.58%  "
$ws.Range("D12:E12").Style = "Normal"
$ws.Range("D13").Value = "'0.07378"
$ws.Range("E13").Value = "'  +4.33%  "
$ws.Range("D13:E13").Style = "Normal"
$ws.Range("D14").Value = "'5.492"
$ws.Range("E14").Value = "'  +1.95%  "
$ws.Range("D14:E14").Style = "Normal"
$ws.Range("D15").Value = "'93.40"
$ws.Range("E15").Value = "'  +1.81%  "
$ws.Range("D15:E15").Style = "Normal"
$ws.Range("D16").Value = "'6.605"
$ws.Range("E16").Value = "'  +1.24%  "
$ws.Range("D16:E16").Style = "Normal"
$ws.Range("D17").Value = "'1.016"
$ws.Range("E17").Value = "'  +1.31%  "
$ws.Range("D17:E17").Style = "Normal"
$ws.Range("D18").Value = "'0.000008857"
$ws.Range("E18").Value = "'  +1.55%  "
$ws.Range("D18:E18").Style = "Normal"
$ws.Range("D19").Value = "'1.014"
$ws.Range("E19").Value = "'  +1.30%  "
$ws.Range("D19:E19").Style = "Normal"
$ws.Range("D20").Value = "'14.86"
$ws.Range("E20").Value = "'  +0.81%  "
$ws.Range("D20:E20").Style = "Normal"
$ws.Range("D21").Value = "'27.444.42"
$ws.Range("E21").Value = "'  +1.80%  "
$ws.Range("D21:E21").Style = "Normal"
$ws.Range("D22").Value = "'5.362"
$ws.Range("E22").Value = "'  +0.74%  "
$ws.Range("D22:E22").Style = "Normal"
$ws.Range("D23").Value = "'10.76"
$ws.Range("E23").Value = "'  +1.39%  "
$ws.Range("D23:E23").Style = "Normal"
$ws.Range("D24").Value = "'2.078.63"
$ws.Range("E24").Value = "'  +0.66%  "
$ws.Range("D24:E24").Style = "Normal"
$ws.Range("D25").Value = "'1.906"
$ws.Range("E25").Value = "'  +0.77%  "
$ws.Range("D25:E25").Style = "Normal"
$ws.Range("D26").Value = "'152.80"
$ws.Range("E26").Value = "'  +1.25%  "
$ws.Range("D26:E26").Style = "Normal"
$ws.Range("E27").Value = "'  +1.28%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'2.173"
$ws.Range("E28").Value = "'  -0.38%  "
$ws.Range("D28:E28").Style = "Normal"
$ws.Range("D29").Value = "'5.284"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Value = "'118.23"
$ws.Range("E30").Value = "'  +1.69%  "
$ws.Range("D30:E30").Style = "Normal"
$ws.Range("D31").Value = "'0.08991"
$ws.Range("E31").Value = "'  +0.33%  "
$ws.Range("D31:E31").Style = "Normal"
$ws.Range("D32").Value = "'0.7609"
$ws.Range("E32").Value = "'  -1.13%  "
$ws.Range("D32:E32").Style = "Normal"
$ws.Range("D33").Value = "'1.179"
$ws.Range("E33").Value = "'  +1.21%  "
$ws.Range("D33:E33").Style = "Normal"
$ws.Range("D34").Value = "'4.570"
$ws.Range("E34").Value = "'  +1.25%  "
$ws.Range("D34:E34").Style = "Normal"
$ws.Range("E35").Value = "'  +1.26%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'1.015"
$ws.Range("E36").Value = "'  +1.37%  "
$ws.Range("D36:E36").Style = "Normal"
$ws.Range("D37").Value = "'1.107"
$ws.Range("E37").Value = "'  +1.89%  "
$ws.Range("D37:E37").Style = "Normal"
$ws.Range("D38").Value = "'0.05376"
$ws.Range("E38").Value = "'  +1.40%  "
$ws.Range("D38:E38").Style = "Normal"
$ws.Range("D39").Value = "'0.01971"
$ws.Range("E39").Value = "'  +0.33%  "
$ws.Range("D39:E39").Style = "Normal"
$ws.Range("D40").Value = "'3.002"
$ws.Range("E40").Value = "'  +1.89%  "
$ws.Range("D40:E40").Style = "Normal"
$ws.Range("D41").Value = "'7.327"
$ws.Range("E41").Value = "'  +0.61%  "
$ws.Range("D41:E41").Style = "Normal"
$ws.Range("D42").Value = "'0.5371"
$ws.Range("E42").Value = "'  +0.63%  "
$ws.Range("D42:E42").Style = "Normal"
$ws.Range("D43").Value = "'2.388"
$ws.Range("E43").Value = "'  +1.45%  "
$ws.Range("D43:E43").Style = "Normal"
$ws.Range("D44").Value = "'0.1670"
$ws.Range("E44").Value = "'  +0.70%  "
$ws.Range("D44:E44").Style = "Normal"
$ws.Range("D45").Value = "'8.574"
$ws.Range("E45").Value = "'  +1.28%  "
$ws.Range("D45:E45").Style = "Normal"
$ws.Range("D46").Value = "'0.4992"
$ws.Range("E46").Value = "'  +1.24%  "
$ws.Range("D46:E46").Style = "Normal"
$ws.Range("D47").Value = "'10.60"
$ws.Range("E47").Value = "'  +0.89%  "
$ws.Range("D47:E47").Style = "Normal"
$ws.Range("D48").Value = "'1.016"
$ws.Range("E48").Value = "'  +1.44%  "
$ws.Range("D48:E48").Style = "Normal"
$ws.Range("D49").Value = "'105.29"
$ws.Range("E49").Value = "'  +1.49%  "
$ws.Range("D49:E49").Style = "Normal"
$ws.Range("D50").Value = "'1.686"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Value = "'0.06329"
$ws.Range("D51").Style = "Normal"
